$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "LookingGlass"
$ws.Range("C5").Value = "/Users/fmdec/Box/Github/Research/"
$ws.Range("B5").Value = "fmdec"

$ws.Range("C8").Select()
